$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Ccl11 -> Ackr2 -> ECs) updated TPM-derived stats
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.259924
$ws.Range("H2").Value = 0.779772
$ws.Range("I2").Value = 0.0006491957374851489
$ws.Range("J2").Value = 0.0006491957374851488
$ws.Range("Q2").Value = 3.211639831810667
$ws.Range("R2").Value = 28.904758486296
$ws.Range("S2").Value = 0.0006491957374851489
$ws.Range("T2").Value = 0.0006491957374851488

# Row 3 (FAPs) updated stats (small floating point recalculation)
$ws.Range("G3").Value = 383.1307676666667
$ws.Range("I3").Value = 0.9569214896224009
$ws.Range("J3").Value = 0.9569214896224006
$ws.Range("Q3").Value = 4733.991606125118
$ws.Range("R3").Value = 42605.92445512606
$ws.Range("S3").Value = 0.9569214896224009
$ws.Range("T3").Value = 0.9569214896224006

# Row 4 (MuSCs) updated stats (small floating point recalculation)
$ws.Range("G4").Value = 16.98778433333333
$ws.Range("H4").Value = 50.963353
$ws.Range("I4").Value = 0.042429314640114
$ws.Range("J4").Value = 0.04242931464011399
$ws.Range("S4").Value = 0.042429314640114
$ws.Range("T4").Value = 0.04242931464011399
